$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new data rows right before the existing row 823 (shifts old rows
# 823..879 down to 825..881, matching the new dimension A1:T881).
$ws.Rows("823:824").Insert()

# New row 823: Plátano "Pintón" quality entry for 2023-01-05 (serial 44931)
$ws.Cells.Item(823, 1).Value = 5
$ws.Cells.Item(823, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(823, 3).Value = "Maule"
$ws.Cells.Item(823, 4).Value = 44931
$ws.Cells.Item(823, 5).Value = 7
$ws.Cells.Item(823, 6).Value = "Fruta"
$ws.Cells.Item(823, 7).Value = 100108
$ws.Cells.Item(823, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(823, 9).Value = 100108006
$ws.Cells.Item(823, 10).Value = "Plátano"
$ws.Cells.Item(823, 11).Value = "Sin especificar"
$ws.Cells.Item(823, 12).Value = "Pintón"
$ws.Cells.Item(823, 13).Value = 800
$ws.Cells.Item(823, 14).Value = 16000
$ws.Cells.Item(823, 15).Value = 16000
$ws.Cells.Item(823, 16).Value = 16000
$ws.Cells.Item(823, 17).Value = "`$/caja 20 kilos"
$ws.Cells.Item(823, 18).Value = "Ecuador"
$ws.Cells.Item(823, 19).Value = 800
$ws.Cells.Item(823, 20).Value = 20

# New row 824: Plátano "Primera Pintón" quality entry for the same date
$ws.Cells.Item(824, 1).Value = 5
$ws.Cells.Item(824, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(824, 3).Value = "Maule"
$ws.Cells.Item(824, 4).Value = 44931
$ws.Cells.Item(824, 5).Value = 7
$ws.Cells.Item(824, 6).Value = "Fruta"
$ws.Cells.Item(824, 7).Value = 100108
$ws.Cells.Item(824, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(824, 9).Value = 100108006
$ws.Cells.Item(824, 10).Value = "Plátano"
$ws.Cells.Item(824, 11).Value = "Sin especificar"
$ws.Cells.Item(824, 12).Value = "Primera Pintón"
$ws.Cells.Item(824, 13).Value = 450
$ws.Cells.Item(824, 14).Value = 17000
$ws.Cells.Item(824, 15).Value = 17000
$ws.Cells.Item(824, 16).Value = 17000
$ws.Cells.Item(824, 17).Value = "`$/caja 20 kilos"
$ws.Cells.Item(824, 18).Value = "Ecuador"
$ws.Cells.Item(824, 19).Value = 850
$ws.Cells.Item(824, 20).Value = 20
